# Task Enumeration spreadsheet update
# Fills in the "Product Definition" task breakdown (rows 6-12) and the
# "Realization" sub-task (row 14) with the real task/hour data, renames the
# old "Flow Diagram" shared string to "Featuere Diagram", and moves the
# active selection to L20 (matching the saved workbook state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: rename "Flow Diagram" -> "Featuere Diagram"; its old numeric
# estimate (G6 = 8) is superseded by the detailed breakdown rows below.
$ws.Range("B6").Value = "Featuere Diagram"
$ws.Range("G6").ClearContents()

# Row 7
$ws.Range("C7").Value = "Familirizing inkscape"
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 6
$ws.Range("H7").Value = "Hours"

# Row 8
$ws.Range("C8").Value = "Preparing version one"
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 4
$ws.Range("H8").Value = "Hours"

# Row 9
$ws.Range("C9").Value = "Revisualizing the whole project"
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = "Days"

# Row 10
$ws.Range("C10").Value = "Feature diagram vesion 2"
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = 3.5
$ws.Range("H10").Value = "Hours"

# Row 11
$ws.Range("B11").Value = "Product Definition"
$ws.Range("C11").Value = "Refering to similar format"
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 4
$ws.Range("H11").Value = "Hours"

# Row 12
$ws.Range("C12").Value = "Preparing version 1 on Notepad++"
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = "Days"

# Row 14 (Realization sub-task gets its own breakdown row too)
$ws.Range("C14").Value = "Major Components BOM"
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = 10
$ws.Range("H14").Value = "Hours"

# Move the active cell/selection to L20
[void]$ws.Range("L20").Select()
